$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 8 already has a date (A8 = 8/31/2025) but B8:D8 are blank.
# Fill them with "-" like the earlier "skip day" row (row 6).
$ws.Range("B8").Value = "-"
$ws.Range("C8").Value = "-"
$ws.Range("D8").Value = "-"

# Row 9: new entry for "Workshop Automation"
# Copy the date cell style from the row above so the number format (style index)
# matches the rest of the date column instead of creating a new custom format.
$ws.Range("A7").Copy() | Out-Null
$ws.Range("A9").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("A9").Value = (Get-Date -Year 2025 -Month 9 -Day 1 -Hour 0 -Minute 0 -Second 0)
$ws.Range("B9").Value = "Workshop Automation"
$ws.Range("C9").Value = "this workflow is designed to send event registration and other event related things automatically "
$ws.Range("D9").Value = "Workshop.json"

# Row 10: new entry, only Name filled in so far (day 10 of n8n)
$ws.Range("A10").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("A10").Value = (Get-Date -Year 2025 -Month 9 -Day 2 -Hour 0 -Minute 0 -Second 0)
$ws.Range("B10").Value = "Automated Shop invoice "
$excel.CutCopyMode = 0

# Match the final selection noted in the diff
$ws.Range("C10").Select()

$wb.Save()
